$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145, shifting existing rows 145:226 down to 146:227
$ws.Rows.Item(145).Insert()

# Populate the newly inserted row 145 with the new record's data
$ws.Cells.Item(145, 1).Value = 10
$ws.Cells.Item(145, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(145, 3).Value = "La Araucanía"
$ws.Cells.Item(145, 4).Value = 44518
$ws.Cells.Item(145, 5).Value = 9
$ws.Cells.Item(145, 6).Value = 100114013
$ws.Cells.Item(145, 7).Value = "Zanahoria"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 280
$ws.Cells.Item(145, 11).Value = 8000
$ws.Cells.Item(145, 12).Value = 8000
$ws.Cells.Item(145, 13).Value = 8000
$ws.Cells.Item(145, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(145, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(145, 16).Value = 400
$ws.Cells.Item(145, 17).Value = 20
$ws.Cells.Item(145, 18).Value = "Hortaliza"

# Ensure the date cell keeps the existing date number format used by column D
$ws.Cells.Item(145, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
